# Add a new "Save" column (H) to the s_vals sheet, matching the header
# styling used by the existing columns (B:G) and filling in the save
# flag values for each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style of an existing header cell (G1) onto the new header
# cell (H1) so it picks up the same bold/centered/bordered formatting.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Set the header text and the data values for the new "Save" column.
$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 0
